$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark (Word auto-creates this to mark the
#    location of the last edit; it is stripped on a clean save). It is
#    a hidden bookmark so it won't show up while iterating
#    $d.Bookmarks, but it can still be addressed directly by name.
# ------------------------------------------------------------------
Try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
    Write-Output "Removed _GoBack bookmark"
} Catch {
    Write-Output "No _GoBack bookmark present"
}

# ------------------------------------------------------------------
# 2) Remove the trailing blank paragraphs at the very end of the
#    document body (after "...Administrative Coordination Unit."):
#      - one paragraph whose pPr only carries an rFonts "Muli" rPr
#      - two completely empty paragraphs
#    Find the last paragraph that still has visible text, then
#    delete everything from the end of that paragraph through the
#    end of the document range. That collapses the trailing empty
#    paragraph marks so <w:sectPr> immediately follows the closing
#    </w:p> of the real content, matching a clean final save.
# ------------------------------------------------------------------
$lastReal = $null
For ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    If ($p.Range.Text.Trim().Length -gt 0) {
        $lastReal = $p
        Break
    }
}

If ($lastReal -ne $null) {
    $finalParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
    If ($finalParagraph.Range.End -gt $lastReal.Range.End) {
        $trailingRange = $d.Range($lastReal.Range.End, $finalParagraph.Range.End)
        $trailingRange.Delete()
        Write-Output "Removed trailing empty paragraphs"
    }
}
